$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header label updates: "As of January 2023" -> "As of April 2023" ---
$ws.Range("A3").Value = "As of April 2023"
$ws.Range("A40").Value = "As of April 2023"
$ws.Range("A77").Value = "As of April 2023"
$ws.Range("A114").Value = "As of April 2023"
$ws.Range("A150").Value = "As of April 2023"
$ws.Range("A187").Value = "As of April 2023"
$ws.Range("A224").Value = "As of April 2023"

# --- Data revisions (Q4 2022 SNA Update) ---
$ws.Range("X12").Value = 357179.8389284021
$ws.Range("W13").Value = 104207.75546279758
$ws.Range("X13").Value = 126689.8479333912
$ws.Range("X14").Value = 98114.21356190383
$ws.Range("X15").Value = 37859.5156958617
$ws.Range("X16").Value = 138718.9321919001
$ws.Range("X17").Value = 37584.062590927904
$ws.Range("X18").Value = 30251.798575179804
$ws.Range("X19").Value = 6474.9008879042985
$ws.Range("W20").Value = 23740.156032172017
$ws.Range("X20").Value = 24496.64732426228
$ws.Range("X21").Value = 10496.642977801548
$ws.Range("X22").Value = 1716.0826899361323
$ws.Range("X23").Value = 1689.1264763330616
$ws.Range("X24").Value = 1407.2953550406928
$ws.Range("W25").Value = 96277.78720832548
$ws.Range("X25").Value = 103352.24337818011
$ws.Range("W26").Value = 267856.9990140938
$ws.Range("X26").Value = 299676.3628743383
$ws.Range("W27").Value = 201537.43066244532
$ws.Range("X27").Value = 236142.04416248604
$ws.Range("X28").Value = 75083.11550096364
$ws.Range("W29").Value = 1383.4624608284107
$ws.Range("X29").Value = 1206.9313653727734
$ws.Range("X30").Value = 269643.5618345118
$ws.Range("X31").Value = 245414.8783651068
$ws.Range("W33").Value = 1954486.610016682
$ws.Range("X33").Value = 2103198.042669804
$ws.Range("W50").Value = 104595.9596207038
$ws.Range("X50").Value = 105247.6872728677
$ws.Range("W57").Value = 19703.198060002815
$ws.Range("X60").Value = 1735.4696454060925
$ws.Range("W62").Value = 92805.21750277742
$ws.Range("X62").Value = 92143.95604362548
$ws.Range("W63").Value = 179028.02112954867
$ws.Range("X63").Value = 183067.35521554653
$ws.Range("W64").Value = 174980.03439215332
$ws.Range("X64").Value = 186754.67682011292
$ws.Range("W66").Value = 2709.553777323661
$ws.Range("X66").Value = 2706.7157979275776
$ws.Range("X67").Value = 215554.7390334859
$ws.Range("W70").Value = 1775210.3171187984
$ws.Range("X70").Value = 1783498.9562229416
$ws.Range("W86").Value = 3.278799167724273
$ws.Range("V87").Value = 20.018223659035115
$ws.Range("W87").Value = 21.574298736930174
$ws.Range("W88").Value = 7.6188453724914496
$ws.Range("W89").Value = 43.54911750489066
$ws.Range("W90").Value = -22.335025099994425
$ws.Range("W91").Value = 0.34694384234698816
$ws.Range("W92").Value = 14.201501560144948
$ws.Range("W93").Value = 14.072210832567222
$ws.Range("V94").Value = -12.212155915657647
$ws.Range("W94").Value = 3.1865472622213815
$ws.Range("W95").Value = -4.579031152338757
$ws.Range("W96").Value = 7.804617373356763
$ws.Range("W97").Value = 12.116893091796726
$ws.Range("W98").Value = -7.546532514907284
$ws.Range("V99").Value = -1.2230608487609231
$ws.Range("W99").Value = 7.34796298812617
$ws.Range("V100").Value = 11.889890721863864
$ws.Range("W100").Value = 11.879235553807675
$ws.Range("V101").Value = 7.116210738549327
$ws.Range("W101").Value = 17.170315899283224
$ws.Range("W102").Value = 13.552733123165254
$ws.Range("V103").Value = -28.402050333728027
$ws.Range("W103").Value = -12.760092915707403
$ws.Range("W104").Value = 8.610683540404722
$ws.Range("W105").Value = 12.648242551455112
$ws.Range("V107").Value = 6.894628059165896
$ws.Range("W107").Value = 7.608720975164559
$ws.Range("V124").Value = 2.029649787629893
$ws.Range("W124").Value = 0.623090657160418
$ws.Range("V131").Value = -1.5802395421134747
$ws.Range("W131").Value = 0.008090300872211742
$ws.Range("W134").Value = 6.85219897322736
$ws.Range("V136").Value = -0.49752927539149994
$ws.Range("W136").Value = -0.712526167111406
$ws.Range("V137").Value = -17.306247613759936
$ws.Range("W137").Value = 2.256258020678729
$ws.Range("V138").Value = -0.30011719262513736
$ws.Range("W138").Value = 6.729134823217066
$ws.Range("V140").Value = 4.50618518365107
$ws.Range("W140").Value = -0.10473973315586704
$ws.Range("W141").Value = -3.567651074529394
$ws.Range("V144").Value = -0.2909612687144687
$ws.Range("W144").Value = 0.4669102598274577
$ws.Range("X159").Value = 94.56544110858061
$ws.Range("W160").Value = 99.62885358161638
$ws.Range("X160").Value = 120.37304687268997
$ws.Range("X161").Value = 116.6396461788216
$ws.Range("X162").Value = 137.49575248590102
$ws.Range("X163").Value = 103.53662044448902
$ws.Range("X164").Value = 113.72382466823183
$ws.Range("X165").Value = 101.62631805048001
$ws.Range("X166").Value = 109.94704283483144
$ws.Range("W167").Value = 120.48884632776526
$ws.Range("X167").Value = 124.31822264345396
$ws.Range("X168").Value = 106.89490909399657
$ws.Range("X169").Value = 84.05270794205455
$ws.Range("X170").Value = 97.32964680795747
$ws.Range("X171").Value = 79.84387303853259
$ws.Range("W172").Value = 103.74178284259085
$ws.Range("X172").Value = 112.16388769899142
$ws.Range("W173").Value = 149.6173600781
$ws.Range("X173").Value = 163.69732469314062
$ws.Range("W174").Value = 115.1773865873025
$ws.Range("X174").Value = 126.44504982862858
$ws.Range("X175").Value = 103.99917370791789
$ws.Range("W176").Value = 51.05868251837814
$ws.Range("X176").Value = 44.5902508973004
$ws.Range("X177").Value = 125.09284789726813
$ws.Range("X178").Value = 123.49842220699425
$ws.Range("W180").Value = 110.09887623844212
$ws.Range("X180").Value = 117.92538679831439
$ws.Range("W196").Value = 17.694694051584108
$ws.Range("X196").Value = 16.98270118561908
$ws.Range("W197").Value = 5.33172010126527
$ws.Range("X197").Value = 6.0236765802886945
$ws.Range("W198").Value = 4.6645623158682845
$ws.Range("X198").Value = 4.665001182549478
$ws.Range("W199").Value = 1.349403485815997
$ws.Range("X199").Value = 1.800092760061851
$ws.Range("W200").Value = 9.138561343024849
$ws.Range("X200").Value = 6.595619118008021
$ws.Range("W201").Value = 1.9163148851182432
$ws.Range("X201").Value = 1.7869958904687178
$ws.Range("W202").Value = 1.3553351022735711
$ws.Range("X202").Value = 1.4383713735668044
$ws.Range("W203").Value = 0.2904160183274038
$ws.Range("X203").Value = 0.30785978098786393
$ws.Range("W204").Value = 1.214649203044138
$ws.Range("X204").Value = 1.164733269396076
$ws.Range("W205").Value = 0.5628256797185255
$ws.Range("X205").Value = 0.49908010395811747
$ws.Range("W206").Value = 0.08144569782100966
$ws.Range("X206").Value = 0.08159396571887893
$ws.Range("W207").Value = 0.07708296571764316
$ws.Range("X207").Value = 0.08031228833728284
$ws.Range("W208").Value = 0.07788061110818496
$ws.Range("X208").Value = 0.06691216549699092
$ws.Range("W209").Value = 4.925988580065213
$ws.Range("X209").Value = 4.91405190007616
$ws.Range("W210").Value = 13.70472417878614
$ws.Range("X210").Value = 14.248604115945662
$ws.Range("W211").Value = 10.31152782677417
$ws.Range("X211").Value = 11.227760742051988
$ws.Range("W212").Value = 3.3830778653294113
$ws.Range("X212").Value = 3.5699498562509593
$ws.Range("W213").Value = 0.07078393137810254
$ws.Range("X213").Value = 0.05738553102877047
$ws.Range("W214").Value = 12.702371461723736
$ws.Range("X214").Value = 12.82064533933408
$ws.Range("W215").Value = 11.14663469525601
$ws.Range("X215").Value = 11.668652850854533
$ws.Range("W233").Value = 21.494091045792292
$ws.Range("X233").Value = 21.17783814802725
$ws.Range("W234").Value = 5.892031981340955
$ws.Range("X234").Value = 5.901191413969714
$ws.Range("W235").Value = 4.669266695583748
$ws.Range("X235").Value = 4.716424362461995
$ws.Range("W236").Value = 1.8806753288397249
$ws.Range("X236").Value = 1.5438777646960857
$ws.Range("W237").Value = 7.591927508150831
$ws.Range("X237").Value = 7.512230284573221
$ws.Range("W238").Value = 1.944180455956722
$ws.Range("X238").Value = 1.853017010334412
$ws.Range("W239").Value = 1.6387113358815055
$ws.Range("X239").Value = 1.6690607693101112
$ws.Range("W240").Value = 0.33719699065071024
$ws.Range("X240").Value = 0.330199729038432
$ws.Range("W241").Value = 1.1099078159922762
$ws.Range("X241").Value = 1.104839004208768
$ws.Range("W242").Value = 0.5687375694071373
$ws.Range("X242").Value = 0.5505801426974535
$ws.Range("W243").Value = 0.11042560885290895
$ws.Range("X243").Value = 0.1144757879120391
$ws.Range("W244").Value = 0.0914921312905452
$ws.Range("X244").Value = 0.09730701772214295
$ws.Range("W245").Value = 0.09999904956622419
$ws.Range("X245").Value = 0.09882590482475047
$ws.Range("W246").Value = 5.227843518474033
$ws.Range("X246").Value = 5.166470982341706
$ws.Range("W247").Value = 10.084890753683467
$ws.Range("X247").Value = 10.264505879119936
$ws.Range("W248").Value = 9.856862181611776
$ws.Range("X248").Value = 10.471252375477597
$ws.Range("W249").Value = 3.9016880641946354
$ws.Range("X249").Value = 4.047990985306547
$ws.Range("W250").Value = 0.1526328317943375
$ws.Range("X250").Value = 0.1517643612004016
$ws.Range("W251").Value = 12.59171880522346
$ws.Range("X251").Value = 12.086059163722945
$ws.Range("W252").Value = 10.755720327712712
$ws.Range("X252").Value = 11.142088913054504
